# "Schulferien - unnötige Zeilen gelöscht"
# The footnote / legend / contact-info block that used to live in rows
# 21:31 (below the actual Schulferien table in A1:G17) is unnecessary and
# gets removed. The cells keep their existing formatting/styles - only
# their contents (the shared-string text) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21:G31").ClearContents()

# Reflect the author's final selection/scroll position in the sheet view.
$ws.Range("A19:I39").Select()
